# "Added a few Testing cases"
#
# Adds several new test-case rows to the "Sheet 1 - Test Cases" sheet:
#   - Managers section: "A manager should be able to delete a user's reservation"
#   - Admins section: two extra rows about the admin/manager assignment page,
#     and one extra row about admins assigning users to themselves
#   - Reservations section: two rows about the 10-reservation limit
#   - Two new rows for a brand-new "Feedback" section at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom of the sheet upward so row numbers used below always
# refer to the same (not-yet-shifted) rows.

# --- New "Feedback" rows appended after the existing "New Users" row (old row 33) ---
$ws.Rows(34).Insert()
$ws.Rows(34).RowHeight = 18
$ws.Range("A34").Value = "Feedback"
$ws.Range("B34").Value = "Anyone should be able to give a rating to a room that they had previously reserved"

$ws.Rows(35).Insert()
$ws.Rows(35).RowHeight = 18
$ws.Range("B35").Value = "Anyone should be able to leave a comment on a room they previously reserved"

# --- Reservations section: two new rows after "Users can delete reservations" (old row 31) ---
$ws.Rows(32).Insert()
$ws.Rows(32).RowHeight = 18
$ws.Range("B32").Value = "Users should not be able to make more than 10 reservations"

$ws.Rows(33).Insert()
$ws.Rows(33).RowHeight = 18
$ws.Range("B33").Value = "A recurring reservation should be limited based on limit of reservations"

# --- Admins section: expand the merged A25:A27 block with 3 new rows ---
# after "Admins can change any user's privilege level" stays last; insert the
# new rows between the existing ones instead.

# New row after "Admins can assign users to different managers" (old row 26)
$ws.Rows(27).Insert()
$ws.Rows(27).RowHeight = 18
$ws.Range("B27").Value = "Admins can assign users to themselves as managers"

# New row after "Admins can assign perform the tasks of managers" (old row 25)
$ws.Rows(26).Insert()
$ws.Rows(26).RowHeight = 18
$ws.Range("B26").Value = "An Admin can view/update/delete managers reservations as well"

$ws.Rows(26).Insert()
$ws.Rows(26).RowHeight = 18
$ws.Range("B26").Value = "Admin manager assign page should display both admins and managers"

# --- Managers section: new row after "A manager can see a list of their users" (old row 22) ---
$ws.Rows(23).Insert()
$ws.Rows(23).RowHeight = 18
$ws.Range("B23").Value = "A manager should be able to delete a user's reservation"
